# Update column F ("dSF") values for several rows to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 3
    4  = -6
    9  = -4
    14 = 2
    15 = -3
    19 = -8
    21 = 8
    28 = -6
    30 = 7
    33 = -1
    34 = -4
    35 = -5
    38 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
